# Apply updated Betfair odds values to Sheet1, per diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA2").Value = 130
$ws.Range("AF2").Value = 11.5
$ws.Range("AI2").Value = 60
$ws.Range("AK2").Value = 16.5
$ws.Range("AN2").Value = 7.4
$ws.Range("AO2").Value = 60
$ws.Range("N2").Value = 5.3
$ws.Range("O2").Value = 1.21
$ws.Range("Q2").Value = 1.65
$ws.Range("T2").Value = 1.71
$ws.Range("U2").Value = 2.32
$ws.Range("AA3").Value = 540
$ws.Range("F3").Value = 1.41
$ws.Range("G3").Value = 1.48
$ws.Range("I3").Value = 11
$ws.Range("P3").Value = 1.86
$ws.Range("R3").Value = 1.32
$ws.Range("V3").Value = 1.1
$ws.Range("W3").Value = 3.05
$ws.Range("AE4").Value = 18
$ws.Range("AK4").Value = 48
$ws.Range("AO4").Value = 10.5
$ws.Range("H4").Value = 1.91
$ws.Range("I4").Value = 1.94
$ws.Range("K4").Value = 4.1
$ws.Range("Q4").Value = 1.76
$ws.Range("S4").Value = 2.9
$ws.Range("U4").Value = 2.34
$ws.Range("V4").Value = 2.06
$ws.Range("F5").Value = 1.66
$ws.Range("H5").Value = 5.7
$ws.Range("K5").Value = 4.5
$ws.Range("U5").Value = 2.42
$ws.Range("Z5").Value = 48
$ws.Range("F6").Value = 1.41
$ws.Range("N6").Value = 7.6
$ws.Range("T6").Value = 1.52
$ws.Range("L9").Value = 1.23
$ws.Range("AB10").Value = 8.4
$ws.Range("AC10").Value = 8.6
$ws.Range("AG10").Value = 13.5
$ws.Range("AJ10").Value = 36
$ws.Range("AN10").Value = 34
$ws.Range("N10").Value = 2.6
$ws.Range("O10").Value = 1.52
$ws.Range("X10").Value = 10.5
$ws.Range("AJ13").Value = 21
$ws.Range("G13").Value = 1.93
$ws.Range("K13").Value = 4
$ws.Range("L13").Value = 1.35
$ws.Range("AJ14").Value = 60
$ws.Range("F14").Value = 3.25
$ws.Range("H14").Value = 2.24
$ws.Range("I14").Value = 2.26
$ws.Range("U14").Value = 2.62
$ws.Range("V14").Value = 1.79
$ws.Range("Z14").Value = 16.5
$ws.Range("AJ15").Value = 38
$ws.Range("AO15").Value = 12.5
$ws.Range("F15").Value = 2.48
$ws.Range("R15").Value = 1.9
$ws.Range("S15").Value = 2.06
$ws.Range("T15").Value = 1.42
$ws.Range("U15").Value = 3.25
$ws.Range("Y15").Value = 23
$ws.Range("AA16").Value = 480
$ws.Range("AC16").Value = 17.5
$ws.Range("AD16").Value = 50
$ws.Range("AE16").Value = 190
$ws.Range("AI16").Value = 1000
$ws.Range("H16").Value = 14.5
$ws.Range("I16").Value = 15
$ws.Range("K16").Value = 8
$ws.Range("S16").Value = 1.96
$ws.Range("T16").Value = 1.86
$ws.Range("U16").Value = 2.1
$ws.Range("Z16").Value = 160
$ws.Range("AB17").Value = 10
$ws.Range("AI17").Value = 150
$ws.Range("G17").Value = 1.34
$ws.Range("J17").Value = 6
$ws.Range("P17").Value = 2.58
$ws.Range("R17").Value = 1.61
$ws.Range("P19").Value = 1.85
$ws.Range("N20").Value = 2.56
$ws.Range("AA21").Value = 40
$ws.Range("AC21").Value = 7.6
$ws.Range("AD21").Value = 12
$ws.Range("AE21").Value = 32
$ws.Range("AF21").Value = 22
$ws.Range("AG21").Value = 14.5
$ws.Range("AH21").Value = 19.5
$ws.Range("AI21").Value = 50
$ws.Range("AJ21").Value = 60
$ws.Range("AK21").Value = 40
$ws.Range("AL21").Value = 60
$ws.Range("AM21").Value = 130
$ws.Range("AN21").Value = 44
$ws.Range("AO21").Value = 30
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 3.35
$ws.Range("I21").Value = 2.76
$ws.Range("J21").Value = 3.25
$ws.Range("K21").Value = 3.4
$ws.Range("M21").Value = 1.09
$ws.Range("N21").Value = 3.1
$ws.Range("O21").Value = 1.4
$ws.Range("R21").Value = 1.26
$ws.Range("T21").Value = 1.87
$ws.Range("U21").Value = 1.96
$ws.Range("V21").Value = 1.57
$ws.Range("W21").Value = 1.43
$ws.Range("X21").Value = 11.5
$ws.Range("Z21").Value = 16
$ws.Range("G22").Value = 2.44
$ws.Range("J22").Value = 2.72
$ws.Range("K22").Value = 3.55
$ws.Range("P22").Value = 1.67
$ws.Range("V22").Value = 1.33
$ws.Range("W22").Value = 1.7
